# feat: add 2022-Q3 data
#
# - Insert a new worksheet "2022-Q3" right after "总计" (becomes the 2nd
#   sheet), carrying the quarterly fund-holding data for 2022-Q3. The
#   existing "2022-Q2" and "2022-Q1" sheets are untouched, just pushed
#   one/two places to the right.
# - Update the "总计" (summary) sheet with a new top data row for 2022-Q3,
#   shifting the 2022-Q2 / 2022-Q1 rows down by one.

$wb = $excel.ActiveWorkbook

# ---- 1. insert the new "2022-Q3" sheet right after "总计" --------------
$wsTotal0 = $wb.Worksheets.Item("总计")
$wsNew = $wb.Worksheets.Add($null, $wsTotal0)
$wsNew.Name = "2022-Q3"

# re-fetch every sheet handle by name now that the sheet collection has
# shifted - stale handles grabbed before Add() can resolve to the wrong
# (shifted) position.
$wsTotal = $wb.Worksheets.Item("总计")
$wsQ3    = $wb.Worksheets.Item("2022-Q3")
$wsQ2    = $wb.Worksheets.Item("2022-Q2")
$wsQ1    = $wb.Worksheets.Item("2022-Q1")

# clone the header row + index cell formatting from the "2022-Q2" sheet
# (same layout/styling is reused by every quarter sheet) - copy the
# header (B1:H1) and the index cell (A2) separately so an empty A1 cell
# isn't materialised.
$wsQ2.Range("B1:H1").Copy($wsQ3.Range("B1:H1"))
$wsQ2.Range("A2").Copy($wsQ3.Range("A2"))

# data row - fund codes / ratios are stored as TEXT (leading zeros /
# fixed decimals must survive), so force the number format to Text
# before writing them.
$wsQ3.Range("B2:G2").NumberFormat = "@"
$wsQ3.Range("B2").Value = "005029"
$wsQ3.Range("C2").Value = "中银产业精选混合"
$wsQ3.Range("D2").Value = "0.44"
$wsQ3.Range("E2").Value = "83.28"
$wsQ3.Range("F2").Value = "2.99"
$wsQ3.Range("G2").Value = "0.0132"
$wsQ3.Range("H2").Value = 10

# ---- 2. update the "总计" summary sheet --------------------------------
# Row 2 used to describe 2022-Q2 and row 3 described 2022-Q1; row 2 now
# becomes 2022-Q3, row 3 becomes 2022-Q2, and a brand new row 4 is added
# for 2022-Q1 - the index column (A) already lines up (0,1,2 by row) so
# it only needs to be written for the newly added row.
$wsTotal.Range("B2").Value = "2022-Q3"
$wsTotal.Range("D2").Value = 0.01

$wsTotal.Range("B3").Value = "2022-Q2"
$wsTotal.Range("D3").Value = 1.17

$wsTotal.Range("A3").Copy($wsTotal.Range("A4"))
$wsTotal.Range("A4").Value = 2
$wsTotal.Range("B4").Value = "2022-Q1"
$wsTotal.Range("C4").Value = 1
$wsTotal.Range("D4").Value = 0.87

# "2022-Q1" was (and remains) the active tab - inserting the new sheet
# shifts the active tab to it, so restore the original selection.
$wsQ1.Activate()
